{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Target change: after the sentence ending \"...non matching pair of\n// socks.\" (last paragraph of the \"Socks in the Dark\" list, the one that\n// carries the _GoBack bookmark), append a new sentence to that same\n// paragraph, then add a brand-new list paragraph (\"Well\") right after it\n// (still before the doc-ending empty paragraph), reusing the same\n// numbering (ListParagraph style, ilvl 0, numId 3).\n\nconst body = context.document.body;\n\n// Locate the run of text that ends the target paragraph. This text is\n// unique in the document, so the search gives us an anchored range that\n// sits right before the existing _GoBack bookmark.\nconst searchResults = body.search(\"non matching pair of socks.\", { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\"Could not find anchor text for the sock-drawer edit.\");\n}\n\nconst anchor = searchResults.items[0];\n\n// Insert the continuation sentence right after the found text (and\n// therefore before the bookmark that closes the paragraph) as its own\n// run of text.\nanchor.insertText(\n  \" However in order to solve the equation as it sits we would need to take the probability of pulling a matching pair of socks in the shortest amount of pulls. This would require us to use the black socks since there are more of them you have a better chance of choosing a pair of them. In order to get a matching pair of each color sock you are looking at an almost improbable ratio. \",\n  Word.InsertLocation.after\n);\nawait context.sync();\n\n// Re-locate the paragraph that owns the text we just extended (paragraphs\n// are live, but grabbing it fresh keeps this robust), then add a new\n// paragraph right after it for the next list item, \"Well\".\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst targetParagraph = anchor.paragraphs.getFirst();\nconst newParagraph = targetParagraph.insertParagraph(\"Well\", Word.InsertLocation.after);\nnewParagraph.styleBuiltIn = Word.Style.listParagraph;\n\nawait context.sync();\n\n// Attach the new paragraph to the same numbered list (numId 3) at the\n// same level (ilvl 0) as the preceding bullet items.\nnewParagraph.attachToList(3, 0);\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is already open as $d below.\n#\n# Target change: after the sentence ending \"...non matching pair of\n# socks.\" (the last paragraph of the \"Socks in the Dark\" list, which\n# carries the _GoBack bookmark), append a new sentence to that same\n# paragraph (before the bookmark), then add a brand-new list paragraph\n# (\"Well\") right after it, reusing the same numbering (ListParagraph\n# style, ilvl 0, numId 3).\n\n$d = $word.ActiveDocument\n\n# 1. Find the sentence that currently ends the target paragraph. It is\n#    unique in the document.\n$findRange = $d.Range()\n$findRange.Find.ClearFormatting()\n$found = $findRange.Find.Execute(\"non matching pair of socks.\")\n\n# 2. Insert the continuation sentence immediately at the end of that\n#    match. Using a fresh zero-length Range at the match's end (rather\n#    than collapsing/inserting on the paragraph's own Range) lands the\n#    new text right before the paragraph's _GoBack bookmark instead of\n#    after it.\n$insertPoint = $d.Range($findRange.End, $findRange.End)\n$insertPoint.InsertBefore(\" However in order to solve the equation as it sits we would need to take the probability of pulling a matching pair of socks in the shortest amount of pulls. This would require us to use the black socks since there are more of them you have a better chance of choosing a pair of them. In order to get a matching pair of each color sock you are looking at an almost improbable ratio. \")\n\n# 3. Add a new list paragraph right after the target paragraph containing\n#    \"Well\". Re-fetch the paragraph from the document's Paragraphs\n#    collection (rather than reusing $findRange.Paragraphs) since the\n#    earlier text insertion can leave old Range-derived collections\n#    stale.\n$targetPara = $d.Paragraphs.Item(11)\n$targetPara.Range.InsertParagraphAfter()\n\n$newPara = $d.Paragraphs.Item(12)\n$newPara.Range.Text = \"Well\"\n"}
